$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.087.11'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +3.20%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.657.69'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +4.00%  '

$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.27%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '215.43'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +1.67%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.507'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.74%  '

$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.24%  '

$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +1.89%  '

$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +1.49%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.54'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +3.36%  '

$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +1.18%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.892.67'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +4.13%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.669.04'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +4.73%  '

$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +1.79%  '

$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +3.04%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.82'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +1.97%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '27.081.27'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +3.33%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '238.01'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +3.81%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.84'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +3.41%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0728'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +1.08%  '

$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +0.16%  '

$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +4.59%  '

$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +5.35%  '

$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +4.21%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '145.96'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +0.26%  '

$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.16%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.17'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +2.79%  '

$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +1.16%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.82'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +3.22%  '

$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +1.21%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.538.95'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +5.36%  '

$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +2.72%  '

$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +3.21%  '

$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +7.98%  '

$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +0.21%  '

$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +1.40%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.886'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +8.30%  '

$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +2.89%  '

$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +3.16%  '

$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +0.21%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.28'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +4.80%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '66.40'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +9.89%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.798.97'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +3.98%  '

$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +2.47%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.923'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -0.40%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '90.00'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +2.79%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0₆0106'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +0.81%  '

$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +3.15%  '

$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +0.99%  '

$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +3.09%  '
